$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated numeric values for existing rows 173-175 ---
$data = @(
    @(173, 8027, 6634, 17546, 10911, 1393, 2364, 971, 31792, 28457, 236, 1264, 5343, 2775, 2713, 1692, 14435, 3334),
    @(174, 7171, 5770, 16759, 10988, 1400, 2364, 963, 31074, 27747, 154, 1461, 4765, 2511, 2759, 1680, 14416, 3327),
    @(175, 6662, 5290, 16504, 11214, 1372, 2332, 960, 31010, 27719, 130, 1427, 4831, 2373, 2711, 1618, 14628, 3292)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value  = $entry[1]
    $ws.Cells.Item($row, 3).Value  = $entry[2]
    $ws.Cells.Item($row, 4).Value  = $entry[3]
    $ws.Cells.Item($row, 5).Value  = $entry[4]
    $ws.Cells.Item($row, 6).Value  = $entry[5]
    $ws.Cells.Item($row, 7).Value  = $entry[6]
    $ws.Cells.Item($row, 8).Value  = $entry[7]
    $ws.Cells.Item($row, 9).Value  = $entry[8]
    $ws.Cells.Item($row, 10).Value = $entry[9]
    $ws.Cells.Item($row, 11).Value = $entry[10]
    $ws.Cells.Item($row, 12).Value = $entry[11]
    $ws.Cells.Item($row, 13).Value = $entry[12]
    $ws.Cells.Item($row, 14).Value = $entry[13]
    $ws.Cells.Item($row, 15).Value = $entry[14]
    $ws.Cells.Item($row, 16).Value = $entry[15]
    $ws.Cells.Item($row, 17).Value = $entry[16]
    $ws.Cells.Item($row, 18).Value = $entry[17]
}

# --- New row 176 ---
# The "Serie" label is a date-like piece of text ("01-07-2021"). Assigning
# it with a plain .Value would make Excel auto-convert it into a date
# serial number. Entering it as a formula that evaluates to text, then
# converting that formula to a static value via copy/paste-special,
# keeps it as genuine text (stored as a shared string) exactly like the
# other "Serie" cells in the column.
$labelCell = $ws.Cells.Item(176, 1)
$labelCell.Formula = '="01-07-2021"'
$labelCell.Copy($labelCell)
$labelCell.PasteSpecial(-4163)  # xlPasteValues

$newRow = @(5797, 4407, 15731, 11324, 1390, 2332, 941, 30327, 27055, 151, 999, 4447, 2502, 2683, 1624, 14648, 3273)
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(176, $i + 2).Value = $newRow[$i]
}
